# Commit: "add app.manifest to grant administrator previleges"
#
# Adds three new rows to the "Issues" sheet documenting the deployment
# issue that required an app.manifest (elevated permissions) fix, and
# marks the corresponding "Deployment" feature row on the "Features"
# sheet as FIXED.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Issues")
$ws2 = $wb.Worksheets.Item("Features")

# --- Issues sheet: row 16 ------------------------------------------------
# Same look as the other resolved ("Fixed") issue rows (e.g. row 13):
# blue fill, centered ID/Status, left-aligned wrapped text columns.
$ws1.Range("A13:F13").Copy()
$ws1.Range("A16:F16").PasteSpecial(-4122) # xlPasteFormats

$ws1.Cells.Item(16, 1).Value = 14
$ws1.Cells.Item(16, 2).Value = "Deployment Issue"
$ws1.Cells.Item(16, 3).Value = "在VM上无法运行"
$ws1.Cells.Item(16, 4).Value = "缺乏权限无法创建DB"
$ws1.Cells.Item(16, 5).Value = "使用app.manifest提升权限"
$ws1.Cells.Item(16, 6).Value = "Fixed"

# --- Issues sheet: row 17 -------------------------------------------------
# Plain style (no fill), matching rows such as row 9 (ID centered, Name
# left-aligned/wrapped).
$ws1.Range("A9:B9").Copy()
$ws1.Range("A17:B17").PasteSpecial(-4122)

$ws1.Cells.Item(17, 1).Value = 15
$ws1.Cells.Item(17, 2).Value = "无需询问权限"

# --- Issues sheet: row 18 -------------------------------------------------
$ws1.Range("A9:C9").Copy()
$ws1.Range("A18:C18").PasteSpecial(-4122)

$ws1.Cells.Item(18, 1).Value = 16
$ws1.Cells.Item(18, 2).Value = "若使用新DB则无法启动程序"
$ws1.Cells.Item(18, 3).Value = "Slice must be between 0.0 and 1.0."

# --- Features sheet: mark "Deployment" row (row 2) as FIXED ---------------
$ws2.Cells.Item(2, 5).Value = "FIXED"

# --- Refresh the saved selections on both sheets --------------------------
$ws1.Activate()
$ws1.Range("C25").Select()
$ws2.Activate()
$ws2.Range("E5").Select()
$ws1.Activate()
